$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Cells.Item(2, 3).Value = "21:19:46"      # C2 Time

$ws.Cells.Item(2, 6).Value = 1.16            # F2 Odd_H_Back
$ws.Cells.Item(2, 7).Value = 1.17            # G2 Odd_H_Lay
$ws.Cells.Item(2, 9).Value = 190             # I2 Odd_A_Lay
$ws.Cells.Item(2, 10).Value = 1.01           # J2 Odd_D_Back
$ws.Cells.Item(2, 11).Value = 9.6            # K2 Odd_D_Lay
$ws.Cells.Item(2, 12).Value = 0              # L2 Odd_Over05_HT_Back
$ws.Cells.Item(2, 13).Value = 0              # M2 Odd_Over05_FT_Back
$ws.Cells.Item(2, 14).Value = 0              # N2 Odd_Under15_FT_Back
$ws.Cells.Item(2, 15).Value = 0              # O2 Odd_Over15_FT_Back
$ws.Cells.Item(2, 16).Value = 2.88           # P2 Odd_Under25_FT_Back
$ws.Cells.Item(2, 17).Value = 1.49           # Q2 Odd_Over25_FT_Back
$ws.Cells.Item(2, 18).Value = 1.43           # R2 Odd_Under35_FT_Back
$ws.Cells.Item(2, 19).Value = 3.15           # S2 Odd_Over35_FT_Back
$ws.Cells.Item(2, 20).Value = 1.88           # T2 Odd_BTTS_Yes_Back
$ws.Cells.Item(2, 21).Value = 1.88           # U2 Odd_BTTS_No_Back
$ws.Cells.Item(2, 23).Value = 6.8            # W2 Double_Chance_Draw_or_Away_Back
$ws.Cells.Item(2, 26).Value = 1000           # Z2 Odd_CS_0x2_Lay
$ws.Cells.Item(2, 29).Value = 1000           # AC2 Odd_CS_1x1_Lay
$ws.Cells.Item(2, 30).Value = 1000           # AD2 Odd_CS_1x2_Lay
$ws.Cells.Item(2, 32).Value = 2.98           # AF2 Odd_CS_2x0_Lay
$ws.Cells.Item(2, 33).Value = 4.4            # AG2 Odd_CS_2x1_Lay
$ws.Cells.Item(2, 34).Value = 14             # AH2 Odd_CS_2x2_Lay
$ws.Cells.Item(2, 35).Value = 65             # AI2 Odd_CS_2x3_Lay
$ws.Cells.Item(2, 36).Value = 10.5           # AJ2 Odd_CS_3x0_Lay
$ws.Cells.Item(2, 37).Value = 16             # AK2 Odd_CS_3x1_Lay
$ws.Cells.Item(2, 38).Value = 50             # AL2 Odd_CS_3x2_Lay
$ws.Cells.Item(2, 39).Value = 210            # AM2 Odd_CS_3x3_Lay
$ws.Cells.Item(2, 40).Value = 32             # AN2 Odd_CS_Goleada_H_Lay
$ws.Cells.Item(2, 41).Value = 380            # AO2 Odd_CS_Goleada_A_Lay

# --- Row 3 updates ---
$ws.Cells.Item(3, 3).Value = "22:10:06"      # C3 Time

$ws.Cells.Item(3, 6).Value = 2.86            # F3 Odd_H_Back
$ws.Cells.Item(3, 7).Value = 3               # G3 Odd_H_Lay
$ws.Cells.Item(3, 8).Value = 3.55            # H3 Odd_A_Back
$ws.Cells.Item(3, 9).Value = 3.8             # I3 Odd_A_Lay
$ws.Cells.Item(3, 10).Value = 2.58           # J3 Odd_D_Back
$ws.Cells.Item(3, 11).Value = 2.64           # K3 Odd_D_Lay
$ws.Cells.Item(3, 12).Value = 2.32           # L3 Odd_Over05_HT_Back
$ws.Cells.Item(3, 13).Value = 1.28           # M3 Odd_Over05_FT_Back
$ws.Cells.Item(3, 14).Value = 1.78           # N3 Odd_Under15_FT_Back
$ws.Cells.Item(3, 15).Value = 2.22           # O3 Odd_Over15_FT_Back
$ws.Cells.Item(3, 16).Value = 1.22           # P3 Odd_Under25_FT_Back
$ws.Cells.Item(3, 17).Value = 5.1            # Q3 Odd_Over25_FT_Back
$ws.Cells.Item(3, 18).Value = 1.06           # R3 Odd_Under35_FT_Back
$ws.Cells.Item(3, 19).Value = 14             # S3 Odd_Over35_FT_Back
$ws.Cells.Item(3, 20).Value = 3.15           # T3 Odd_BTTS_Yes_Back
$ws.Cells.Item(3, 21).Value = 1.4            # U3 Odd_BTTS_No_Back
$ws.Cells.Item(3, 22).Value = 1.36           # V3 Double_Chance_Home_or_Draw_Back
$ws.Cells.Item(3, 23).Value = 1.51           # W3 Double_Chance_Draw_or_Away_Back
$ws.Cells.Item(3, 24).Value = 4.6            # X3 Odd_CS_0x0_Lay
$ws.Cells.Item(3, 25).Value = 7.8            # Y3 Odd_CS_0x1_Lay
$ws.Cells.Item(3, 26).Value = 27             # Z3 Odd_CS_0x2_Lay
$ws.Cells.Item(3, 27).Value = 130            # AA3 Odd_CS_0x3_Lay
$ws.Cells.Item(3, 28).Value = 5.8            # AB3 Odd_CS_1x0_Lay
$ws.Cells.Item(3, 29).Value = 7.6            # AC3 Odd_CS_1x1_Lay
$ws.Cells.Item(3, 30).Value = 25             # AD3 Odd_CS_1x2_Lay
$ws.Cells.Item(3, 31).Value = 130            # AE3 Odd_CS_1x3_Lay
$ws.Cells.Item(3, 32).Value = 17             # AF3 Odd_CS_2x0_Lay
$ws.Cells.Item(3, 33).Value = 20             # AG3 Odd_CS_2x1_Lay
$ws.Cells.Item(3, 34).Value = 55             # AH3 Odd_CS_2x2_Lay
$ws.Cells.Item(3, 35).Value = 320            # AI3 Odd_CS_2x3_Lay
$ws.Cells.Item(3, 36).Value = 75             # AJ3 Odd_CS_3x0_Lay
$ws.Cells.Item(3, 37).Value = 90             # AK3 Odd_CS_3x1_Lay
$ws.Cells.Item(3, 38).Value = 250            # AL3 Odd_CS_3x2_Lay
$ws.Cells.Item(3, 40).Value = 180            # AN3 Odd_CS_Goleada_H_Lay
$ws.Cells.Item(3, 41).Value = 310            # AO3 Odd_CS_Goleada_A_Lay
